$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Row 8 (panel renamed FIRECLASS 64-4 -> FC64-4, trip current updated)
$ws.Cells.Item(8, 1).Value = "FC64-4"
$ws.Cells.Item(8, 6).Value = 0.337

# Row 9 (panel renamed FIRECLASS 32-1 -> FC32-1, trip current & volt drop updated)
$ws.Cells.Item(9, 1).Value = "FC32-1"
$ws.Cells.Item(9, 6).Value = 0.198
$ws.Cells.Item(9, 7).Value = 0.388

# Row 8 volt drop now stored as text "0.530"
$ws.Cells.Item(8, 7).Value = "'0.530"

# Row 4: NGC-1928 -> NGC-1928/T958 OR TC71687
$ws.Cells.Item(4, 2).Value = "NGC-1928/T958 OR TC71687"

# Selected cell moved from G11 to G10
$ws.Range("G10").Select()
